$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ICON0007-001")

# Set RUN column (A) value for rows 3 and 4 to "run", matching row 2's existing value
$ws.Range("A3").Value = "run"
$ws.Range("A4").Value = "run"

# Update the active cell selection to A4
$ws.Activate()
$ws.Range("A4").Select()
